$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for each record.
# Update the date serial from 45186 (2023-09-17) to 45188 (2023-09-19)
# for rows 2 through 12, keeping the existing date formatting.
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value = 45188
    }
}
